# Apply crypto price-sheet refresh (GitHub Actions daily update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as plain text in the sheet
# (e.g. "3.924.24", "  +1.59%  "). Force text format before writing so
# Excel does not reinterpret numeric-looking strings as numbers, then
# clear the temporary format back to the default (unstyled) cell style.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "68.183.94"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").Value = "3.924.81"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D5").Value = "485.03"
$ws.Range("E5").Value = "  +4.35%  "

$ws.Range("D6").Value = "145.97"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("E10").Value = "  +3.97%  "

$ws.Range("D11").Value = "0.0000362"
$ws.Range("E11").Value = "  +6.88%  "

$ws.Range("D12").Value = "42.60"
$ws.Range("E12").Value = "  -0.81%  "

$ws.Range("D13").Value = "10.59"
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").Value = "4.548.96"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "14.86"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").Value = "3.950.25"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "19.90"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").Value = "68.316.23"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("D21").Value = "446.52"
$ws.Range("E21").Value = "  +3.59%  "

$ws.Range("D22").Value = "14.89"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("D23").Value = "3.38"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("D24").Value = "88.65"
$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").Value = "11.41"
$ws.Range("E25").Value = "  +12.94%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "10.74"
$ws.Range("E26").Value = "  +11.79%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("D28").Value = "38.83"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  +3.04%  "

$ws.Range("D30").Value = "694.73"
$ws.Range("E30").Value = "  -6.01%  "

$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("E33").Value = "  +2.87%  "

$ws.Range("D34").Value = "0.0₃0947"
$ws.Range("E34").Value = "  +20.49%  "

$ws.Range("D35").Value = "41.65"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("D36").Value = "59.05"
$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("E37").Value = "  -4.82%  "

$ws.Range("D38").Value = "5.65"
$ws.Range("E38").Value = "  +5.03%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "0.0479"
$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  +15.87%  "

$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("D43").Value = "0.369"
$ws.Range("E43").Value = "  +10.30%  "

$ws.Range("D44").Value = "2.97"
$ws.Range("E44").Value = "  +6.63%  "

$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("E48").Value = "  -1.18%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.18"
$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "145.32"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  -1.66%  "

$numRng.ClearFormats()
